# Refresh the "cryptos" price/volume snapshot (GitHub Actions scrape update).
# Columns: A=index, B=Coin, C=Link, D=Price, E=Volume(1h); also re-sorts the
# Chainlink/Uniswap rows (19-20) to match their new rank order.
#
# Every Price/Volume cell in this sheet is plain text (prices use '.' as a
# thousands separator, e.g. "72.373.33", and volumes keep padding spaces,
# e.g. "  -0.11%  "). For Price values that *do* look like a clean Excel
# number (e.g. "584.03"), writing via Range.Value would otherwise make COM
# silently coerce the cell to a Number, so such values are entered with a
# leading apostrophe (forcing text entry, same as typing it in the UI) and
# the cell's style is reset to "Normal" right after so the quote-prefix
# formatting COM applies doesn't leave a stray style index behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '72.373.33'
$ws.Range("E2").Value = '  -0.11%  '
$ws.Range("D3").Value = '2.641.80'
$ws.Range("E3").Value = '  -1.35%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = "'584.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.93%  '
$ws.Range("D6").Value = "'175.15"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.53%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E8").Value = '  -0.96%  '
$ws.Range("D9").Value = '2.641.29'
$ws.Range("E9").Value = '  -1.36%  '
$ws.Range("D10").Value = "'0.172"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.52%  '
$ws.Range("E11").Value = '  +1.08%  '
$ws.Range("D12").Value = "'0.357"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.37%  '
$ws.Range("E13").Value = '  -2.70%  '
$ws.Range("D14").Value = '3.125.01'
$ws.Range("E14").Value = '  -0.46%  '
$ws.Range("D15").Value = "'0.0000186"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.29%  '
$ws.Range("D16").Value = '72.280.51'
$ws.Range("E16").Value = '  -0.19%  '
$ws.Range("D17").Value = "'25.85"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.86%  '
$ws.Range("D18").Value = '2.644.66'
$ws.Range("E18").Value = '  -1.35%  '
$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D19").Value = "'12.10"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.43%  '
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").Value = "'7.92"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.35%  '
$ws.Range("D21").Value = "'373.89"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.26%  '
$ws.Range("E22").Value = '  -1.52%  '
$ws.Range("D23").Value = "'2.04"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.20%  '
$ws.Range("E24").Value = '  -0.01%  '
$ws.Range("E25").Value = '  -2.10%  '
$ws.Range("E26").Value = '  -2.21%  '
$ws.Range("D27").Value = "'9.49"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.41%  '
$ws.Range("D28").Value = '2.777.85'
$ws.Range("E28").Value = '  -1.41%  '
$ws.Range("D29").Value = "'0.996"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.24%  '
$ws.Range("E30").Value = '  +0.58%  '
$ws.Range("D31").Value = "'7.96"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.64%  '
$ws.Range("D32").Value = "'495.29"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.36%  '
$ws.Range("E33").Value = '  -2.87%  '
$ws.Range("E34").Value = '  -1.76%  '
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("D36").Value = "'162.65"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.63%  '
$ws.Range("D37").Value = "'19.18"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.68%  '
$ws.Range("E38").Value = '  +4.43%  '
$ws.Range("D39").Value = "'18.85"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.43%  '
$ws.Range("E40").Value = '  -1.86%  '
$ws.Range("E41").Value = '  -0.07%  '
$ws.Range("D42").Value = "'1.72"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -6.44%  '
$ws.Range("E43").Value = '  -0.97%  '
$ws.Range("E44").Value = '  -2.77%  '
$ws.Range("E45").Value = '  -2.47%  '
$ws.Range("D46").Value = "'39.04"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.47%  '
$ws.Range("D47").Value = "'151.47"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.50%  '
$ws.Range("D48").Value = "'3.65"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.27%  '
$ws.Range("E49").Value = '  -1.01%  '
$ws.Range("E50").Value = '  -3.56%  '
$ws.Range("D51").Value = "'0.602"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.64%  '
